# Update column G ("K") values on rows 2-32 of Sheet1 per the regenerated
# save_data output (K column replaces old Strike# values; std/mean and
# s_vals were recalculated upstream, the observable result in this sheet
# is the refreshed K values below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 0
    4  = 0
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 2
    14 = 2
    15 = 0
    16 = 2
    17 = 1
    18 = 0
    19 = 0
    20 = 0
    21 = 2
    22 = 4
    23 = 2
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 2
    31 = 2
    32 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
